# Replace the yearly header labels (currently shared-string text like
# "1960 [YR1960]") in row 1, columns E:BL, with plain numeric year values
# (1960 .. 2019), left-aligned. Column BM ("2020 [YR2020]") is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 5   # column E
$endCol   = 64  # column BL
$startYear = 1960

for ($col = $startCol; $col -le $endCol; $col++) {
    $year = $startYear + ($col - $startCol)
    $ws.Cells.Item(1, $col).Value = $year
}

# Left-align the newly numeric header cells (they were text before, which
# is left-aligned by default display-wise, but Excel now stamps an explicit
# left-alignment style on them since they are numeric).
$headerRange = $ws.Range($ws.Cells.Item(1, $startCol), $ws.Cells.Item(1, $endCol))
$headerRange.HorizontalAlignment = -4131  # xlLeft

# Reflect the selection state recorded after making this edit.
$ws.Range("E1:BL1").Select() | Out-Null
